$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header I1 text and add new header J1, copying the formatting of I1.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "無擔保金額"
$ws.Range("J1").Value = "無擔保資產分類"

$ws.Range("J3").Select() | Out-Null
